# Update daily COVID-19 figures for Valais for the rows dated
# 2021-05-04 .. 2021-05-07 (rows 434-437), correcting the "Nb nouveaux cas
# positifs" (C), and filling in the rest of the data for 2021-05-07 (row 437),
# which had previously been left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 434 (2021-05-04): correct "Nb nouveaux cas positifs" ---
$ws.Range("C434").Value = 87

# --- Row 435 (2021-05-05): correct "Nb nouveaux cas positifs" ---
$ws.Range("C435").Value = 60

# Column L on row 435 was previously stored as text ("1"); normalize it to a
# real number. The column is formatted as Text (@), so we briefly switch to
# a general number format to store a numeric value, then restore the
# original Text format.
$ws.Range("L435").NumberFormat = "general"
$ws.Range("L435").Value = 1
$ws.Range("L435").NumberFormat = "@"

# --- Row 436 (2021-05-06): correct "Nb nouveaux cas positifs" ---
$ws.Range("C436").Value = 37

$ws.Range("L436").NumberFormat = "general"
$ws.Range("L436").Value = 0
$ws.Range("L436").NumberFormat = "@"

$ws.Range("M436").NumberFormat = "general"
$ws.Range("M436").Value = 0
$ws.Range("M436").NumberFormat = "@"

# --- Row 437 (2021-05-07): fill in the previously missing daily data ---
$ws.Range("C437").Value = 5
$ws.Range("E437").Value = 7
$ws.Range("F437").Value = 5
$ws.Range("G437").Value = 23

$ws.Range("L437").NumberFormat = "general"
$ws.Range("L437").Value = 0
$ws.Range("L437").NumberFormat = "@"

$ws.Range("M437").NumberFormat = "general"
$ws.Range("M437").Value = 0
$ws.Range("M437").NumberFormat = "@"

# Restore the selection on the visible (bottom-right) pane to A2, as saved
# in the workbook.
$ws.Range("A2").Select() | Out-Null
